$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1000093.1
$ws.Range("I11").Value = 1000093.1
$ws.Range("K11").Value = 1000093.1
$ws.Range("M11").Value = -999953.1

$ws.Range("H15").Value = 1373.1613
$ws.Range("I15").Value = 1373.1613
$ws.Range("K15").Value = 4119.4839
$ws.Range("M15").Value = -3950.4839

$ws.Range("H33").Value = 50011176
$ws.Range("I33").Value = 90910980
$ws.Range("J33").Value = 22532
$ws.Range("K33").Value = 90910980
$ws.Range("L33").Value = 22532
$ws.Range("M33").Value = -90910751
$ws.Range("N33").Value = -22990

$ws.Range("H40").Value = 1510.8462
$ws.Range("J40").Value = 1300
$ws.Range("L40").Value = 1300
$ws.Range("N40").Value = -1650

$ws.Range("H42").Value = 240
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 240
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 720
$ws.Range("N42").Value = -1180
$ws.Range("M42").ClearContents()

$ws.Range("H51").Value = 385163.94
$ws.Range("I51").Value = 530840.5600000001
$ws.Range("J51").Value = 2762.75
$ws.Range("K51").Value = 530840.5600000001
$ws.Range("L51").Value = 2762.75
$ws.Range("M51").Value = -530356.5600000001
$ws.Range("N51").Value = -3730.75

$ws.Range("H132").Value = 1783928.9
$ws.Range("I132").Value = 2218248
$ws.Range("J132").Value = 3220.2
$ws.Range("K132").Value = 6654744
$ws.Range("L132").Value = 9660.599999999999
$ws.Range("M132").Value = -6652214
$ws.Range("N132").Value = -14720.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6859.3164
$ws.Range("I32").Value = 3364.2024
$ws.Range("K32").Value = 3364.2024
$ws.Range("M32").Value = -3077.2024

$ws.Range("H132").Value = 1561.9803
$ws.Range("I132").Value = 1062.1794
$ws.Range("J132").Value = 3186.3333
$ws.Range("K132").Value = 3186.5382
$ws.Range("L132").Value = 9558.999899999999
$ws.Range("M132").Value = -656.5382
$ws.Range("N132").Value = -14618.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1753.1111
$ws.Range("I94").Value = 1196.3334
$ws.Range("J94").Value = 2866.6667
$ws.Range("K94").Value = 1196.3334
$ws.Range("L94").Value = 2866.6667
$ws.Range("M94").Value = -745.3334
$ws.Range("N94").Value = -3768.6667

$ws.Range("H107").Value = 819.9
$ws.Range("I107").Value = 736.0526
$ws.Range("J107").Value = 2413
$ws.Range("K107").Value = 736.0526
$ws.Range("L107").Value = 2413
$ws.Range("M107").Value = 1183.9474
$ws.Range("N107").Value = -6253

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1900.6492
$ws.Range("I31").Value = 1082.0227
$ws.Range("J31").Value = 4671.385
$ws.Range("K31").Value = 1082.0227
$ws.Range("L31").Value = 4671.385
$ws.Range("M31").Value = -787.0227
$ws.Range("N31").Value = -5261.385

$ws.Range("H34").Value = 1900.6492
$ws.Range("I34").Value = 1082.0227
$ws.Range("J34").Value = 4671.385
$ws.Range("K34").Value = 1082.0227
$ws.Range("L34").Value = 4671.385
$ws.Range("M34").Value = -880.0227
$ws.Range("N34").Value = -5075.385

$ws.Range("H64").Value = 15000
$ws.Range("J64").Value = 15000
$ws.Range("L64").Value = 15000
$ws.Range("N64").Value = -15496

$ws.Range("H67").Value = 15000
$ws.Range("J67").Value = 15000
$ws.Range("L67").Value = 15000
$ws.Range("N67").Value = -16716

$ws.Range("H132").Value = 1126.75
$ws.Range("I132").Value = 774.7143
$ws.Range("J132").Value = 2495.7778
$ws.Range("K132").Value = 2324.1429
$ws.Range("L132").Value = 7487.3334
$ws.Range("M132").Value = 205.8571000000002
$ws.Range("N132").Value = -12547.3334

$ws.Range("H134").Value = 1126.96
$ws.Range("I134").Value = 976.0909
$ws.Range("J134").Value = 2233.3333
$ws.Range("K134").Value = 2928.2727
$ws.Range("L134").Value = 6699.999899999999
$ws.Range("M134").Value = -393.2727
$ws.Range("N134").Value = -11769.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1110.7
$ws.Range("I114").Value = 429.75
$ws.Range("J114").Value = 1564.6666
$ws.Range("K114").Value = 1289.25
$ws.Range("L114").Value = 4693.9998
$ws.Range("M114").Value = 1964.75
$ws.Range("N114").Value = -11201.9998

$ws.Range("H117").Value = 818.7857
$ws.Range("J117").Value = 1041.7778
$ws.Range("L117").Value = 3125.3334
$ws.Range("N117").Value = -10009.3334

$ws.Range("H120").Value = 15377.2
$ws.Range("I120").Value = 1900
$ws.Range("J120").Value = 18746.5
$ws.Range("K120").Value = 5700
$ws.Range("L120").Value = 56239.5
$ws.Range("M120").Value = -862
$ws.Range("N120").Value = -65915.5

$ws.Range("H121").Value = 46879076
$ws.Range("I121").Value = 628.75
$ws.Range("J121").Value = 62505224
$ws.Range("K121").Value = 1886.25
$ws.Range("L121").Value = 187515672
$ws.Range("M121").Value = -576.25
$ws.Range("N121").Value = -187518292

$ws.Range("H129").Value = 1084.875
$ws.Range("J129").Value = 1634.875
$ws.Range("L129").Value = 4904.625
$ws.Range("N129").Value = -14904.625

$ws.Range("H131").Value = 943.5685999999999
$ws.Range("J131").Value = 994.13336
$ws.Range("L131").Value = 2982.40008
$ws.Range("N131").Value = -13062.40008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 7000
$ws.Range("J35").Value = 7000
$ws.Range("L35").Value = 7000
$ws.Range("M35").Value = -7596

$ws.Range("H80").Value = 2254.2307
$ws.Range("I80").Value = 2233.889
$ws.Range("J80").Value = 2300
$ws.Range("K80").Value = 2233.889
$ws.Range("L80").Value = 2300
$ws.Range("M80").Value = -1235.889
$ws.Range("N80").Value = -4296

$ws.Range("H83").Value = 2254.2307
$ws.Range("I83").Value = 2233.889
$ws.Range("J83").Value = 2300
$ws.Range("K83").Value = 11169.445
$ws.Range("L83").Value = 11500
$ws.Range("M83").Value = -6177.445
$ws.Range("N83").Value = -21484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 45000
$ws.Range("I64").Value = 50000
$ws.Range("J64").Value = 40000
$ws.Range("K64").Value = 50000
$ws.Range("L64").Value = 40000
$ws.Range("M64").Value = -49775
$ws.Range("N64").Value = -40450

$ws.Range("H67").Value = 45000
$ws.Range("I67").Value = 50000
$ws.Range("J67").Value = 40000
$ws.Range("K67").Value = 50000
$ws.Range("L67").Value = 40000
$ws.Range("M67").Value = -49220
$ws.Range("N67").Value = -41560

$ws.Range("H68").Value = 11917.091
$ws.Range("I68").Value = 23097.6
$ws.Range("J68").Value = 2600
$ws.Range("K68").Value = 23097.6
$ws.Range("L68").Value = 2600
$ws.Range("M68").Value = -22348.6
$ws.Range("N68").Value = -4098

$ws.Range("H71").Value = 11917.091
$ws.Range("I71").Value = 23097.6
$ws.Range("J71").Value = 2600
$ws.Range("K71").Value = 115488
$ws.Range("L71").Value = 13000
$ws.Range("M71").Value = -111744
$ws.Range("N71").Value = -20488

$ws.Range("H82").Value = 1378.16
$ws.Range("I82").Value = 936.8
$ws.Range("J82").Value = 2040.2
$ws.Range("K82").Value = 936.8
$ws.Range("L82").Value = 2040.2
$ws.Range("M82").Value = -575.8
$ws.Range("N82").Value = -2762.2

$ws.Range("H85").Value = 1378.16
$ws.Range("I85").Value = 936.8
$ws.Range("J85").Value = 2040.2
$ws.Range("K85").Value = 936.8
$ws.Range("L85").Value = 2040.2
$ws.Range("M85").Value = 311.2
$ws.Range("N85").Value = -4536.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 22194.334
$ws.Range("I63").Value = 6500
$ws.Range("J63").Value = 24156.125
$ws.Range("K63").Value = 6500
$ws.Range("L63").Value = 24156.125
$ws.Range("M63").Value = -5876
$ws.Range("N63").Value = -25404.125

$ws.Range("H66").Value = 22194.334
$ws.Range("I66").Value = 6500
$ws.Range("J66").Value = 24156.125
$ws.Range("K66").Value = 19500
$ws.Range("L66").Value = 72468.375
$ws.Range("M66").Value = -16380
$ws.Range("N66").Value = -78708.375
